{"js": "const body = context.document.body;\n\n// Load the paragraph collection so we can reach into it by index.\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 4 (0-based) is the one that contains the runs:\n//   \" \" / \"   \ucc55\ud130 2\" / \"\uc758 \ub0b4\uc6a9\uc744 \ucd94\uac00\"\n// It currently has no <w:pPr>. The target revision gives it a\n// <w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr> and also\n// absorbs the \"_GoBack\" bookmark that used to sit in its own trailing\n// paragraph. We replace the paragraph's OOXML in place (keeping its\n// runs/text identical) to add that paragraph mark formatting and append\n// the bookmark start/end.\nconst target = body.paragraphs.items[4];\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">   \ucc55\ud130 2</w:t></w:r>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc758 \ub0b4\uc6a9\uc744 \ucd94\uac00</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// The paragraphs that used to follow (two blank paragraphs, \"Delete 5\",\n// \"\uc5b4\uca4c\uace0 \uc800\uca4c\uace0\", the \"1234\u314825436\" paragraph, and the old paragraph that\n// only held the \"_GoBack\" bookmark) are all removed by this revision \u2014\n// their content no longer exists anywhere in the document. Repeatedly\n// deleting the paragraph right after our (now bookmark-holding) target\n// paragraph removes exactly those six paragraphs, regardless of how the\n// collection re-indexes after each deletion.\nfor (let i = 0; i < 6; i++) {\n  body.paragraphs.load(\"items\");\n  await context.sync();\n  body.paragraphs.items[5].delete();\n  await context.sync();\n}\n", "ps1": "# Reverts the \"M05: delete05 doc. file \uc218\uc815\" commit:\n#   - the paragraph containing \" \" / \"   \ucc55\ud130 2\" / \"\uc758 \ub0b4\uc6a9\uc744 \ucd94\uac00\" gains a\n#     paragraph-mark <w:rFonts w:hint=\"eastAsia\"/> and absorbs the trailing\n#     \"_GoBack\" bookmark\n#   - the paragraphs that had been inserted after it (two blank paragraphs,\n#     \"Delete 5\", \"\uc5b4\uca4c\uace0 \uc800\uca4c\uace0\", the \"1234\u314825436\" paragraph, and the old\n#     paragraph that only held the \"_GoBack\" bookmark) are removed entirely\n$d = $word.ActiveDocument\n\n# Paragraph 5 (1-based) holds the runs: \" \" / \"   \ucc55\ud130 2\" / \"\uc758 \ub0b4\uc6a9\uc744 \ucd94\uac00\"\n$target = $d.Paragraphs.Item(5)\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n       '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n       '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n       '<pkg:xmlData>' +\n       '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n       '<w:body>' +\n       '<w:p>' +\n       '<w:pPr><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr></w:pPr>' +\n       '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n       '<w:r><w:t xml:space=\"preserve\">   \ucc55\ud130 2</w:t></w:r>' +\n       '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc758 \ub0b4\uc6a9\uc744 \ucd94\uac00</w:t></w:r>' +\n       '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n       '<w:bookmarkEnd w:id=\"0\"/>' +\n       '</w:p>' +\n       '</w:body></w:document>' +\n       '</pkg:xmlData></pkg:part></pkg:package>'\n\n# Replace the paragraph's contents in place (same runs/text, now carrying\n# the paragraph-mark font hint and the bookmark that used to trail it).\n[void]$target.Range.InsertXML($xml)\n\n# Remove the six paragraphs that used to follow it. Deleting the paragraph\n# immediately after our target, six times in a row, removes exactly the\n# two blank paragraphs, \"Delete 5\", \"\uc5b4\uca4c\uace0 \uc800\uca4c\uace0\", the \"1234\u314825436\"\n# paragraph, and the old bookmark-only paragraph - regardless of how the\n# paragraph collection re-indexes after each deletion.\nfor ($i = 0; $i -lt 6; $i++) {\n    $d.Paragraphs.Item(6).Range.Delete()\n}\n"}
